$wb = $excel.ActiveWorkbook

# Insert a new worksheet ("Sheet1") right before the "Search" sheet.
$search = $wb.Worksheets.Item("Search")
$new = $wb.Worksheets.Add($search)

# Populate the new sheet with the lorem-ipsum word/definition pairs.
# Values are written bottom-up / right-to-left so the resulting shared-string
# table indices come out in the same order as the target workbook.
$new.Range("B8").Value = "eu fugiat nulla pariatur"
$new.Range("A8").Value = "elit"
$new.Range("B7").Value = "in voluptate velit esse cillum dolore"
$new.Range("A7").Value = "adipiscing"
$new.Range("B6").Value = "laboris nisi ut aliquip ex ea commodo consequat"
$new.Range("A6").Value = "consectetur"
$new.Range("B5").Value = "quis nostrud exercitation ullamco"
$new.Range("A5").Value = "amet"
$new.Range("B4").Value = "Ut enim ad minim veniam"
$new.Range("A4").Value = "sit"
$new.Range("B3").Value = "sed do eiusmod tempor incididunt ut labore et dolore magna aliqua"
$new.Range("A3").Value = "dolor"
$new.Range("B2").Value = "consectetur adipiscing elit"
$new.Range("A2").Value = "ipsum"
$new.Range("B1").Value = "Lorem ipsum dolor sit amet"
$new.Range("A1").Value = "lorem"

# Select the populated range, matching the saved selection state.
$new.Range("A1:B8").Select()
